# Commit: "finished upload post process and converted csv to JSOn"
#
# Source OOXML diff shows:
#   1. xl/workbook.xml  : <sheet name="Sheet1" .../>  -> <sheet name="Summary" .../>
#   2. xl/workbook.xml  : <calcPr calcId="162913"/>    -> <calcPr calcId="162913" calcOnSave="0"/>
#   3. xl/workbook.xml  : x15ac:absPath url tweaked (author's local folder path - machine/
#                         save-location metadata Excel stamps on save; not meaningful user
#                         content, left untouched)
#   4. xl/styles.xml    : custom numFmtId's renumbered 168-171 -> 164-167 (Excel's own internal
#                         bookkeeping when it rewrites styles.xml; same 4 format codes / same
#                         cell formatting either way)
#   5. xl/worksheets/sheet1.xml : stored cursor/selection (<selection activeCell="E7" .../>)
#                         removed from the sheet view.
#
# The only real, user-visible content edit here is the worksheet rename (#1) - that is what we
# reproduce below. We also nudge the calculation/selection state in the direction the diff wants
# wherever the object model exposes a hook for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename "Sheet1" -> "Summary"
$ws.Name = "Summary"

# 2) The workbook no longer wants to force a recalculation on every save.
$excel.CalculateBeforeSave = $false
$wb.UpdateRemoteReferences = $false

# 5) Collapse the saved selection back down to the top-left cell (closest available
#    approximation to "no stored selection" for this sheet).
[void]$ws.Range("A1").Select()
